$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the date stored as text in A2 into a real date serial value,
# formatted with the built-in date number format (numFmtId 14).
$ws.Range("A2").Value = 45852
$ws.Range("A2").NumberFormat = "mm-dd-yy"

# Widen column A to fit the "Fecha consulta"-length content.
$ws.Range("A1").ColumnWidth = 23.25

# Move/confirm the active selection to A3, like the saved workbook.
$ws.Range("A3").Select() | Out-Null
